$wb = $excel.ActiveWorkbook

# --- "Nodos Loads" sheet: add a new row of load data (row 6) ---
$loads = $wb.Worksheets.Item("Nodos Loads")
$loads.Cells.Item(6, 1).Value = 9
$loads.Cells.Item(6, 2).Value = 0
$loads.Cells.Item(6, 3).Value = 0
$loads.Cells.Item(6, 4).Value = -10

# Select the new cell on "Nodos Loads" and make it the active sheet/tab
$loads.Activate()
$loads.Range("C7").Select()

# --- "Props" sheet: selection moved (no longer the active tab) ---
$props = $wb.Worksheets.Item("Props")
$props.Activate()
$props.Range("E24").Select()

# Re-activate "Nodos Loads" last so it ends up as the selected/active tab
$loads.Activate()
